# Add an "Email Address" field to the Company record (sheet "Company").
# Inserts a new row above the current row 7 ("Contact Person's Name"),
# shifting the existing rows 7-8 down to 8-9, and fills in the new row
# with the Email Address field definition (Text, required).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Company")

$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "Email Address"
$ws.Range("B7").Value = "Text"
$ws.Range("C7").Value = "Yes"

# The row-insert carries the neighboring row's formatting into the blank
# D7/E7 cells; clear them so no stray formatted-but-empty cells are left
# behind (matching the source record's layout, which only populates A-C).
$ws.Range("D7:E7").Clear()
